$d = $word.ActiveDocument

# The document currently ends with an empty paragraph right before the
# section properties. Add two new paragraphs of text after it.

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter("Merge a branch back to the master branch within the Browser Git Repo view by going to the branch and then adding a “Pull Request” and inserting a message for the pull through to master – This would usually be done by the Repo Admin, or someone who is in control of managing any merge conflicts")

$d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter("You can see the repo tree (branches etc) within the Browser by going to the repo, choosing “Insights…Network”")
